# Add a new "unseen class" evaluation row (row 11) plus a "비고" (Remarks)
# column (J) to the model-summary sheet, per commit:
#   "unseen class 추가 14개 class_EffNetV2M_Test acc 84.08%"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 11 data: unseen-class evaluation of EffNetV2M ---------------
# Set plain values first (no explicit style touches) so the new row
# naturally inherits the same per-column formatting as the rows above it
# (A:border/center, B-E:border/center, F-G:percent, H-I:0.0000_ numfmt).
$ws.Cells.Item(11, 1).Value = 220620
$ws.Cells.Item(11, 2).Value = "박영서"
$ws.Cells.Item(11, 3).Value = "EffNetV2M"
$ws.Cells.Item(11, 4).Value = 14
$ws.Cells.Item(11, 5).Value = 5
$ws.Cells.Item(11, 6).Value = 0.8714
$ws.Cells.Item(11, 7).Value = 0.8408
$ws.Cells.Item(11, 8).Value = 0.5629
$ws.Cells.Item(11, 9).Value = 0.5932

# --- New column J: 비고 (Remarks) header + remark on the new row ---------
# Copy an existing formatted neighbour cell first (carries the style index
# over verbatim, without minting a duplicate style record), then overwrite
# its value with the real text.
$ws.Cells.Item(1, 9).Copy($ws.Cells.Item(1, 10))
$ws.Cells.Item(1, 10).Value = "비고"

$ws.Cells.Item(11, 9).Copy($ws.Cells.Item(11, 10))
$ws.Cells.Item(11, 10).Value = "unseen('맛집' 검색) class acc 51%로 낮음"

# --- Widen column J to fit the remark text --------------------------------
# (39.33203125 isn't exactly reachable through the ColumnWidth->pixel
# round-trip; 38.55 lands on the nearest attainable stored width.)
$ws.Columns.Item(10).ColumnWidth = 38.55

# --- Update the active selection to reflect the new last-used cell -------
$ws.Range("F12").Select()
